$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto snapshot values.
# Cells whose new text would otherwise be auto-parsed as a number are forced to
# the Text format first, so values like "0.1890" keep their trailing zero.

$ws.Range("D2").Value = "27.859.26"
$ws.Range("E2").Value = "  +0.61%  "

$ws.Range("D3").Value = "1.871.24"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  +0.52%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.56"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("E6").Value = "  +0.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4754"
$ws.Range("E7").Value = "  +0.66%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3923"
$ws.Range("E8").Value = "  -0.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.78"
$ws.Range("E9").Value = "  -2.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07973"
$ws.Range("E10").Value = "  -0.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.008"
$ws.Range("E11").Value = "  -1.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.67"
$ws.Range("E12").Value = "  -1.71%  "

$ws.Range("D13").Value = "1.879.81"
$ws.Range("E13").Value = "  -0.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.010"
$ws.Range("E14").Value = "  +0.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.162"
$ws.Range("E15").Value = "  +0.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.011"
$ws.Range("E16").Value = "  +0.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.25"
$ws.Range("E17").Value = "  +1.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06681"
$ws.Range("E18").Value = "  +0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001043"
$ws.Range("E19").Value = "  -0.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.99"
$ws.Range("E20").Value = "  -1.10%  "

$ws.Range("D22").Value = "27.858.44"
$ws.Range("E22").Value = "  +0.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.489"
$ws.Range("E23").Value = "  -0.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.94"
$ws.Range("E24").Value = "  -0.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.328"
$ws.Range("E25").Value = "  +0.95%  "

$ws.Range("D26").Value = "2.099.67"
$ws.Range("E26").Value = "  -0.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.82"
$ws.Range("E27").Value = "  -0.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.72"
$ws.Range("E28").Value = "  -2.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.091"
$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.432"
$ws.Range("E30").Value = "  -2.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.23"
$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9704"
$ws.Range("E32").Value = "  -0.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09481"
$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("E34").Value = "  +0.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.310"
$ws.Range("E35").Value = "  -0.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.345"
$ws.Range("E36").Value = "  -7.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06043"
$ws.Range("E37").Value = "  -0.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02225"
$ws.Range("E38").Value = "  -1.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.201"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.149"
$ws.Range("E40").Value = "  -1.32%  "

$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5940"
$ws.Range("E42").Value = "  -1.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1890"
$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.30"
$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5649"
$ws.Range("E46").Value = "  -1.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.06"
$ws.Range("E47").Value = "  -0.75%  "

$ws.Range("E48").Value = "  -1.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.298"
$ws.Range("E49").Value = "  -2.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06775"
$ws.Range("E50").Value = "  -1.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.97"
$ws.Range("E51").Value = "  -3.22%  "
